# Apply updated cryptocurrency market data scraped on Fri Sep 20 19:07:25 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    if ($val -match "^-?[0-9]*\.?[0-9]+$") {
        $ws.Range($addr).Value = "'" + $val
    } else {
        $ws.Range($addr).Value = $val
    }
}

Set-TextCell 'D2' '62.918.71'
Set-TextCell 'E2' '  -0.98%  '
Set-TextCell 'D3' '2.546.28'
Set-TextCell 'E3' '  +2.99%  '
Set-TextCell 'E4' '  -0.06%  '
Set-TextCell 'D5' '566.98'
Set-TextCell 'E5' '  -0.21%  '
Set-TextCell 'D6' '146.37'
Set-TextCell 'E6' '  +1.90%  '
Set-TextCell 'E7' '  -0.01%  '
Set-TextCell 'D8' '0.583'
Set-TextCell 'E8' '  -1.54%  '
Set-TextCell 'D9' '2.545.42'
Set-TextCell 'E9' '  +3.03%  '
Set-TextCell 'D10' '0.105'
Set-TextCell 'E10' '  -1.39%  '
Set-TextCell 'E11' '  -2.35%  '
Set-TextCell 'E12' '  +0.61%  '
Set-TextCell 'D13' '0.353'
Set-TextCell 'E13' '  -0.63%  '
Set-TextCell 'D14' '27.23'
Set-TextCell 'E14' '  +3.16%  '
Set-TextCell 'D15' '2.998.82'
Set-TextCell 'E15' '  +2.80%  '
Set-TextCell 'D16' '62.882.18'
Set-TextCell 'E16' '  -0.86%  '
Set-TextCell 'D17' '0.0000142'
Set-TextCell 'E17' '  -1.14%  '
Set-TextCell 'D18' '2.538.51'
Set-TextCell 'E18' '  +2.55%  '
Set-TextCell 'D19' '11.47'
Set-TextCell 'E19' '  +1.66%  '
Set-TextCell 'D20' '334.58'
Set-TextCell 'E20' '  -2.12%  '
Set-TextCell 'D21' '4.28'
Set-TextCell 'E21' '  -0.69%  '
Set-TextCell 'D22' '6.75'
Set-TextCell 'E22' '  -1.18%  '
Set-TextCell 'E23' '  +0.04%  '
Set-TextCell 'D24' '64.74'
Set-TextCell 'E24' '  -1.61%  '
Set-TextCell 'E25' '  -3.56%  '
Set-TextCell 'E26' '  +4.83%  '
Set-TextCell 'B27' 'SuiNetwork'
Set-TextCell 'C27' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell 'D27' '1.49'
Set-TextCell 'E27' '  +11.92%  '
Set-TextCell 'B28' 'Binance-PegBSC-USD'
Set-TextCell 'C28' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 'D28' '0.999'
Set-TextCell 'E28' '  -0.12%  '
Set-TextCell 'D29' '8.34'
Set-TextCell 'E29' '  +1.88%  '
Set-TextCell 'D30' '7.27'
Set-TextCell 'E30' '  +6.34%  '
Set-TextCell 'D31' '0.0₃0808'
Set-TextCell 'E31' '  -1.97%  '
Set-TextCell 'E32' '  -0.75%  '
Set-TextCell 'D33' '176.89'
Set-TextCell 'E33' '  -0.16%  '
Set-TextCell 'D34' '1.58'
Set-TextCell 'E34' '  +4.28%  '
Set-TextCell 'D35' '407.89'
Set-TextCell 'E35' '  +9.37%  '
Set-TextCell 'D36' '0.399'
Set-TextCell 'E36' '  -0.44%  '
Set-TextCell 'D37' '19.01'
Set-TextCell 'E37' '  +0.16%  '
Set-TextCell 'E38' '  -0.03%  '
Set-TextCell 'D39' '4.36'
Set-TextCell 'E39' '  -2.41%  '
Set-TextCell 'E40' '  +1.94%  '
Set-TextCell 'D41' '1.00'
Set-TextCell 'E41' '  +0.00%  '
Set-TextCell 'D42' '39.06'
Set-TextCell 'E42' '  -3.46%  '
Set-TextCell 'D43' '153.22'
Set-TextCell 'E43' '  +1.86%  '
Set-TextCell 'D44' '3.74'
Set-TextCell 'E44' '  +0.56%  '
Set-TextCell 'D45' '20.71'
Set-TextCell 'E45' '  -0.51%  '
Set-TextCell 'D46' '0.605'
Set-TextCell 'E46' '  +0.56%  '
Set-TextCell 'D47' '0.0957'
Set-TextCell 'E47' '  -0.94%  '
Set-TextCell 'D48' '0.0518'
Set-TextCell 'E48' '  -1.12%  '
Set-TextCell 'D49' '0.0236'
Set-TextCell 'E49' '  +4.15%  '
Set-TextCell 'D50' '18.24'
Set-TextCell 'E50' '  +0.41%  '
Set-TextCell 'E51' '  -1.00%  '
